$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = 45946.01041666666
$ws.Cells.Item(2, 2).Value2 = 82
$ws.Cells.Item(3, 1).Value2 = 45946.02083333334
$ws.Cells.Item(3, 2).Value2 = 81
$ws.Cells.Item(4, 1).Value2 = 45946.03125
$ws.Cells.Item(4, 2).Value2 = 81
$ws.Cells.Item(5, 1).Value2 = 45946.04166666666
$ws.Cells.Item(5, 2).Value2 = 80
$ws.Cells.Item(6, 1).Value2 = 45946.05208333334
$ws.Cells.Item(6, 2).Value2 = 91
$ws.Cells.Item(7, 1).Value2 = 45946.0625
$ws.Cells.Item(7, 2).Value2 = 91
$ws.Cells.Item(8, 1).Value2 = 45946.07291666666
$ws.Cells.Item(8, 2).Value2 = 86
$ws.Cells.Item(9, 1).Value2 = 45946.08333333334
$ws.Cells.Item(9, 2).Value2 = 86
$ws.Cells.Item(10, 1).Value2 = 45946.09375
$ws.Cells.Item(10, 2).Value2 = 82
$ws.Cells.Item(11, 1).Value2 = 45946.10416666666
$ws.Cells.Item(11, 2).Value2 = 82
$ws.Cells.Item(12, 1).Value2 = 45946.11458333334
$ws.Cells.Item(12, 2).Value2 = 89
$ws.Cells.Item(13, 1).Value2 = 45946.125
$ws.Cells.Item(13, 2).Value2 = 88
$ws.Cells.Item(14, 1).Value2 = 45946.13541666666
$ws.Cells.Item(14, 2).Value2 = 94
$ws.Cells.Item(15, 1).Value2 = 45946.14583333334
$ws.Cells.Item(15, 2).Value2 = 89
$ws.Cells.Item(16, 1).Value2 = 45946.15625
$ws.Cells.Item(16, 2).Value2 = 90
$ws.Cells.Item(17, 1).Value2 = 45946.16666666666
$ws.Cells.Item(17, 2).Value2 = 92
$ws.Cells.Item(18, 1).Value2 = 45946.17708333334
$ws.Cells.Item(18, 2).Value2 = 98
$ws.Cells.Item(19, 1).Value2 = 45946.1875
$ws.Cells.Item(19, 2).Value2 = 98
$ws.Cells.Item(20, 1).Value2 = 45946.19791666666
$ws.Cells.Item(20, 2).Value2 = 99
$ws.Cells.Item(21, 1).Value2 = 45946.20833333334
$ws.Cells.Item(21, 2).Value2 = 100
$ws.Cells.Item(22, 1).Value2 = 45946.21875
$ws.Cells.Item(22, 2).Value2 = 97
$ws.Cells.Item(23, 1).Value2 = 45946.22916666666
$ws.Cells.Item(23, 2).Value2 = 97
$ws.Cells.Item(24, 1).Value2 = 45946.23958333334
$ws.Cells.Item(24, 2).Value2 = 98
$ws.Cells.Item(25, 1).Value2 = 45946.25
$ws.Cells.Item(25, 2).Value2 = 99
$ws.Cells.Item(26, 1).Value2 = 45946.26041666666
$ws.Cells.Item(26, 2).Value2 = 101
$ws.Cells.Item(27, 1).Value2 = 45946.27083333334
$ws.Cells.Item(27, 2).Value2 = 102
$ws.Cells.Item(28, 1).Value2 = 45946.28125
$ws.Cells.Item(28, 2).Value2 = 103
$ws.Cells.Item(29, 1).Value2 = 45946.29166666666
$ws.Cells.Item(29, 2).Value2 = 104
$ws.Cells.Item(30, 1).Value2 = 45946.30208333334
$ws.Cells.Item(30, 2).Value2 = 111
$ws.Cells.Item(31, 1).Value2 = 45946.3125
$ws.Cells.Item(31, 2).Value2 = 115
$ws.Cells.Item(32, 1).Value2 = 45946.32291666666
$ws.Cells.Item(32, 2).Value2 = 116
$ws.Cells.Item(33, 1).Value2 = 45946.33333333334
$ws.Cells.Item(33, 2).Value2 = 118
$ws.Cells.Item(34, 1).Value2 = 45946.34375
$ws.Cells.Item(34, 2).Value2 = 115
$ws.Cells.Item(35, 1).Value2 = 45946.35416666666
$ws.Cells.Item(35, 2).Value2 = 114
$ws.Cells.Item(36, 1).Value2 = 45946.36458333334
$ws.Cells.Item(36, 2).Value2 = 110
$ws.Cells.Item(37, 1).Value2 = 45946.375
$ws.Cells.Item(37, 2).Value2 = 110
$ws.Cells.Item(38, 1).Value2 = 45946.38541666666
$ws.Cells.Item(38, 2).Value2 = 119
$ws.Cells.Item(39, 1).Value2 = 45946.39583333334
$ws.Cells.Item(39, 2).Value2 = 114
$ws.Cells.Item(40, 1).Value2 = 45946.40625
$ws.Cells.Item(40, 2).Value2 = 111
$ws.Cells.Item(41, 1).Value2 = 45946.41666666666
$ws.Cells.Item(41, 2).Value2 = 109
$ws.Cells.Item(42, 1).Value2 = 45946.42708333334
$ws.Cells.Item(42, 2).Value2 = 99
$ws.Cells.Item(43, 1).Value2 = 45946.4375
$ws.Cells.Item(43, 2).Value2 = 99
$ws.Cells.Item(44, 1).Value2 = 45946.44791666666
$ws.Cells.Item(44, 2).Value2 = 99
$ws.Cells.Item(45, 1).Value2 = 45946.45833333334
$ws.Cells.Item(45, 2).Value2 = 100
$ws.Cells.Item(46, 1).Value2 = 45946.46875
$ws.Cells.Item(46, 2).Value2 = 99
$ws.Cells.Item(47, 1).Value2 = 45946.47916666666
$ws.Cells.Item(47, 2).Value2 = 100
$ws.Cells.Item(48, 1).Value2 = 45946.48958333334
$ws.Cells.Item(48, 2).Value2 = 102
$ws.Cells.Item(49, 1).Value2 = 45946.5
$ws.Cells.Item(49, 2).Value2 = 103
$ws.Cells.Item(50, 1).Value2 = 45946.51041666666
$ws.Cells.Item(50, 2).Value2 = 102
$ws.Cells.Item(51, 1).Value2 = 45946.52083333334
$ws.Cells.Item(51, 2).Value2 = 104
$ws.Cells.Item(52, 1).Value2 = 45946.53125
$ws.Cells.Item(52, 2).Value2 = 106
$ws.Cells.Item(53, 1).Value2 = 45946.54166666666
$ws.Cells.Item(53, 2).Value2 = 107
$ws.Cells.Item(54, 1).Value2 = 45946.55208333334
$ws.Cells.Item(54, 2).Value2 = 108
$ws.Cells.Item(55, 1).Value2 = 45946.5625
$ws.Cells.Item(55, 2).Value2 = 110
$ws.Cells.Item(56, 1).Value2 = 45946.57291666666
$ws.Cells.Item(56, 2).Value2 = 112
$ws.Cells.Item(57, 1).Value2 = 45946.58333333334
$ws.Cells.Item(57, 2).Value2 = 114
$ws.Cells.Item(58, 1).Value2 = 45946.59375
$ws.Cells.Item(58, 2).Value2 = 125
$ws.Cells.Item(59, 1).Value2 = 45946.60416666666
$ws.Cells.Item(59, 2).Value2 = 127
$ws.Cells.Item(60, 1).Value2 = 45946.61458333334
$ws.Cells.Item(60, 2).Value2 = 129
$ws.Cells.Item(61, 1).Value2 = 45946.625
$ws.Cells.Item(61, 2).Value2 = 131
$ws.Cells.Item(62, 1).Value2 = 45946.63541666666
$ws.Cells.Item(62, 2).Value2 = 146
$ws.Cells.Item(63, 1).Value2 = 45946.64583333334
$ws.Cells.Item(63, 2).Value2 = 149
$ws.Cells.Item(64, 1).Value2 = 45946.65625
$ws.Cells.Item(64, 2).Value2 = 153
$ws.Cells.Item(65, 1).Value2 = 45946.66666666666
$ws.Cells.Item(65, 2).Value2 = 156
$ws.Cells.Item(66, 1).Value2 = 45946.67708333334
$ws.Cells.Item(66, 2).Value2 = 183
$ws.Cells.Item(67, 1).Value2 = 45946.6875
$ws.Cells.Item(67, 2).Value2 = 192
$ws.Cells.Item(68, 1).Value2 = 45946.69791666666
$ws.Cells.Item(68, 2).Value2 = 205
$ws.Cells.Item(69, 1).Value2 = 45946.70833333334
$ws.Cells.Item(69, 2).Value2 = 215
$ws.Cells.Item(70, 1).Value2 = 45946.71875
$ws.Cells.Item(70, 2).Value2 = 266
$ws.Cells.Item(71, 1).Value2 = 45946.72916666666
$ws.Cells.Item(71, 2).Value2 = 284
$ws.Cells.Item(72, 1).Value2 = 45946.73958333334
$ws.Cells.Item(72, 2).Value2 = 302
$ws.Cells.Item(73, 1).Value2 = 45946.75
$ws.Cells.Item(73, 2).Value2 = 320
$ws.Cells.Item(74, 1).Value2 = 45946.76041666666
$ws.Cells.Item(74, 2).Value2 = 377
$ws.Cells.Item(75, 1).Value2 = 45946.77083333334
$ws.Cells.Item(75, 2).Value2 = 392
$ws.Cells.Item(76, 1).Value2 = 45946.78125
$ws.Cells.Item(76, 2).Value2 = 407
$ws.Cells.Item(77, 1).Value2 = 45946.79166666666
$ws.Cells.Item(77, 2).Value2 = 423
$ws.Cells.Item(78, 1).Value2 = 45946.80208333334
$ws.Cells.Item(78, 2).Value2 = 464
$ws.Cells.Item(79, 1).Value2 = 45946.8125
$ws.Cells.Item(79, 2).Value2 = 474
$ws.Cells.Item(80, 1).Value2 = 45946.82291666666
$ws.Cells.Item(80, 2).Value2 = 484
$ws.Cells.Item(81, 1).Value2 = 45946.83333333334
$ws.Cells.Item(81, 2).Value2 = 494
$ws.Cells.Item(82, 1).Value2 = 45946.84375
$ws.Cells.Item(82, 2).Value2 = 520
$ws.Cells.Item(83, 1).Value2 = 45946.85416666666
$ws.Cells.Item(83, 2).Value2 = 526
$ws.Cells.Item(84, 1).Value2 = 45946.86458333334
$ws.Cells.Item(84, 2).Value2 = 531
$ws.Cells.Item(85, 1).Value2 = 45946.875
$ws.Cells.Item(85, 2).Value2 = 535
$ws.Cells.Item(86, 1).Value2 = 45946.88541666666
$ws.Cells.Item(86, 2).Value2 = 546
$ws.Cells.Item(87, 1).Value2 = 45946.89583333334
$ws.Cells.Item(87, 2).Value2 = 548
$ws.Cells.Item(88, 1).Value2 = 45946.90625
$ws.Cells.Item(88, 2).Value2 = 551
$ws.Cells.Item(89, 1).Value2 = 45946.91666666666
$ws.Cells.Item(89, 2).Value2 = 554
$ws.Cells.Item(90, 1).Value2 = 45946.92708333334
$ws.Cells.Item(90, 2).Value2 = 555
$ws.Cells.Item(91, 1).Value2 = 45946.9375
$ws.Cells.Item(91, 2).Value2 = 558
$ws.Cells.Item(92, 1).Value2 = 45946.94791666666
$ws.Cells.Item(92, 2).Value2 = 558
$ws.Cells.Item(93, 1).Value2 = 45946.95833333334
$ws.Cells.Item(93, 2).Value2 = 556
$ws.Cells.Item(94, 1).Value2 = 45946.96875
$ws.Cells.Item(95, 1).Value2 = 45946.97916666666
$ws.Cells.Item(96, 1).Value2 = 45946.98958333334
$ws.Cells.Item(97, 1).Value2 = 45947
